$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (2024-08-27) entirely; all rows below shift up by one.
$ws.Rows(2).Delete()

# After the shift, the last data row is now row 33 (was 34). Append a new
# row 34 for 2024-09-29, copying the other column values from row 33 and
# formatting to match the rest of the A column.
$last = 33
$new = $last + 1

$ws.Range("A$new").Style = $ws.Range("A$last").Style
$ws.Range("A$new").Value = "'2024-09-29"

$ws.Range("B$new").Value = $ws.Range("B$last").Value2
$ws.Range("C$new").Value = $ws.Range("C$last").Value2
$ws.Range("D$new").Value = $ws.Range("D$last").Value2
$ws.Range("E$new").Value = $ws.Range("E$last").Value2
$ws.Range("F$new").Value = $ws.Range("F$last").Value2
$ws.Range("H$new").Value = $ws.Range("H$last").Value2
$ws.Range("I$new").Value = $ws.Range("I$last").Value2
$ws.Range("J$new").Value = $ws.Range("J$last").Value2
$ws.Range("K$new").Value = $ws.Range("K$last").Value2
